# "Fixed comments in MCADAM script"
# Rewrites the sample-input rows on the Active sheet (rows 2-11) with the
# corrected MCADAM sample set, removes a stray selection on Figure 5,
# and drops the now-unused last row on the Libarkin sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Active": replace rows 2-11 with the corrected sample data.
# ---------------------------------------------------------------------
$active = $wb.Worksheets.Item("Active")

# Clear the old L-column helper values (rows 3-9) - not present any more.
$active.Range("L3:L9").ClearContents()

# --- Row 2 --------------------------------------------------------------
$active.Range("A2").Value = "21-10"
$active.Range("B2").Value = 90
$active.Range("C2").Value = 0
$active.Range("D2").Value = 0
$active.Range("E2").Value = 5
$active.Range("F2").Value = 2.89
$active.Range("G2").Value = 1
$active.Range("H2").Value = 0
$active.Range("I2").ClearContents()
$active.Range("J2").Formula = "=65.91*10^7"
$active.Range("K2").Formula = "=1.2*10^7"
$active.Range("M2").Value = 1
$active.Range("N2").Value = 1
$active.Range("O2").Value = 0
$active.Range("P2").Value = 0.25
$active.Range("Q2").Value = 1

# --- Row 3 --------------------------------------------------------------
$active.Range("A3").Value = "Day-4-02"
$active.Range("B3").Value = 0
$active.Range("C3").Value = 0
$active.Range("D3").Value = 0
$active.Range("E3").Value = 5
$active.Range("F3").Value = 2.89
$active.Range("G3").Value = 1
$active.Range("H3").Value = 0
$active.Range("I3").ClearContents()
$active.Range("J3").Formula = "=74.48*10^7"
$active.Range("K3").Formula = "=2.33*10^7"

# --- Row 4 --------------------------------------------------------------
$active.Range("A4").Value = "Day-04-03"
$active.Range("B4").Value = 90
$active.Range("C4").Value = 0
$active.Range("D4").Value = 1000
$active.Range("E4").Value = 5
$active.Range("F4").Value = 2.89
$active.Range("G4").Value = 1
$active.Range("H4").Value = 0
$active.Range("I4").ClearContents()
$active.Range("J4").Formula = "=91.1*10^7"
$active.Range("K4").Formula = "=1.3*10^7"

# --- Row 5 --------------------------------------------------------------
$active.Range("A5").Value = "Day-4-04"
$active.Range("B5").Value = 0
$active.Range("C5").Value = 0
$active.Range("D5").Value = 1000
$active.Range("E5").Value = 5
$active.Range("F5").Value = 2.89
$active.Range("G5").Value = 1
$active.Range("H5").Value = 0
$active.Range("I5").ClearContents()
$active.Range("J5").Formula = "=167.12*10^7"
$active.Range("K5").Formula = "=3.65*10^7"

# --- Row 6 --------------------------------------------------------------
$active.Range("A6").Value = "21-11"
$active.Range("B6").Value = 90
$active.Range("C6").Value = 0
$active.Range("D6").Value = 2000
$active.Range("E6").Value = 5
$active.Range("F6").Value = 2.89
$active.Range("G6").Value = 1
$active.Range("H6").Value = 0
$active.Range("I6:L6").Clear()

# --- Row 7 --------------------------------------------------------------
$active.Range("A7").Value = "Day-4-03"
$active.Range("B7").Value = 0
$active.Range("C7").Value = 0
$active.Range("D7").Value = 2000
$active.Range("E7").Value = 5
$active.Range("F7").Value = 2.89
$active.Range("G7").Value = 1
$active.Range("H7").Value = 0
$active.Range("I7:L7").Clear()

# --- Row 8 --------------------------------------------------------------
$active.Range("A8").Value = "Day-04-04"
$active.Range("B8").Value = 90
$active.Range("C8").Value = 0
$active.Range("D8").Value = 4000
$active.Range("E8").Value = 5
$active.Range("F8").Value = 2.89
$active.Range("G8").Value = 1
$active.Range("H8").Value = 0
$active.Range("I8:L8").Clear()

# --- Row 9 --------------------------------------------------------------
$active.Range("A9").Value = "Day-4-05"
$active.Range("B9").Value = 0
$active.Range("C9").Value = 0
$active.Range("D9").Value = 4000
$active.Range("E9").Value = 5
$active.Range("F9").Value = 2.89
$active.Range("G9").Value = 1
$active.Range("H9").Value = 0
$active.Range("I9:L9").Clear()

# --- Row 10 -------------------------------------------------------------
$active.Range("A10").Value = "21-12"
$active.Range("B10").Value = 90
$active.Range("C10").Value = 0
$active.Range("D10").Value = 6000
$active.Range("E10").Value = 5
$active.Range("F10").Value = 2.89
$active.Range("G10").Value = 1
$active.Range("H10").Value = 0
$active.Range("I10:L10").Clear()

# --- Row 11 -------------------------------------------------------------
$active.Range("A11").Value = "Day-4-04"
$active.Range("B11").Value = 0
$active.Range("C11").Value = 0
$active.Range("D11").Value = 6000
$active.Range("E11").Value = 5
$active.Range("F11").Value = 2.89
$active.Range("G11").Value = 1
$active.Range("H11").Value = 0
$active.Range("I11:L11").Clear()

# --- Reapply the correct cell formatting (style indices) for rows 2-11 --
# NOTE: multi-area (union) ranges only paste into their first area here,
# so every contiguous block gets its own PasteSpecial call.

# Style "1" (centered, general number format)
$active.Range("A14").Copy() | Out-Null
$active.Range("A2:A11").PasteSpecial(-4122) | Out-Null
$active.Range("D6:D11").PasteSpecial(-4122) | Out-Null
$active.Range("E2:E11").PasteSpecial(-4122) | Out-Null
$active.Range("F2:F11").PasteSpecial(-4122) | Out-Null
$active.Range("G2:G11").PasteSpecial(-4122) | Out-Null
$active.Range("I2:I5").PasteSpecial(-4122) | Out-Null
$active.Range("J2").PasteSpecial(-4122) | Out-Null
$active.Range("K2").PasteSpecial(-4122) | Out-Null

# Style "2" (centered, 2-decimal number format) for J3:K5
$active.Range("J14").Copy() | Out-Null
$active.Range("J3:J5").PasteSpecial(-4122) | Out-Null
$active.Range("K3:K5").PasteSpecial(-4122) | Out-Null

# Style "9" ("Bad"-based highlight, no number format) for B/C/D columns
$active.Range("B19").Copy() | Out-Null
$active.Range("B2").PasteSpecial(-4122) | Out-Null
$active.Range("B6").PasteSpecial(-4122) | Out-Null
$active.Range("B10").PasteSpecial(-4122) | Out-Null
$active.Range("C2").PasteSpecial(-4122) | Out-Null
$active.Range("C6").PasteSpecial(-4122) | Out-Null
$active.Range("C10").PasteSpecial(-4122) | Out-Null
$active.Range("D2").PasteSpecial(-4122) | Out-Null

# Style "10" ("Bad"-based highlight, centered) for B/C/D columns
$active.Range("B14").Copy() | Out-Null
$active.Range("B3:B5").PasteSpecial(-4122) | Out-Null
$active.Range("B7:B9").PasteSpecial(-4122) | Out-Null
$active.Range("B11").PasteSpecial(-4122) | Out-Null
$active.Range("C3:C5").PasteSpecial(-4122) | Out-Null
$active.Range("C7:C9").PasteSpecial(-4122) | Out-Null
$active.Range("C11").PasteSpecial(-4122) | Out-Null
$active.Range("D3:D5").PasteSpecial(-4122) | Out-Null

# Style "13" (erosion column number format) for column H
$active.Range("H14").Copy() | Out-Null
$active.Range("H2:H11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the saved selection on the Active sheet.
$active.Range("P8").Select()

# ---------------------------------------------------------------------
# Sheet "Figure 5": just the saved cursor position moved.
# ---------------------------------------------------------------------
$fig5 = $wb.Worksheets.Item("Figure 5")
$fig5.Range("D18").Select()

# ---------------------------------------------------------------------
# Sheet "Libarkin": drop the now-empty trailing row 11 and move the
# saved selection to the new last row.
# ---------------------------------------------------------------------
$libarkin = $wb.Worksheets.Item("Libarkin")
$libarkin.Range("B11:C11").ClearContents()
$libarkin.Range("B11:C11").Select()
